{"js": "// Fix the typo'd daily-standup date 30/02/2025 -> 30/03/2025 in the\n// \"ATA 30-03-2025\" minutes document (see commit \"Atualizando ata de daily\").\n// The participant name runs (\"Ana Karoline, Lays Abreu, Vitor Restini\" and\n// \"Ana Beatriz Zinatto, Luiz Felipe\") keep the exact same text in the\n// source diff - only Word's live spell-checker re-flagged the proper\n// nouns (wrapping them in <w:proofErr> and re-splitting the runs around\n// them) while resaving, no characters actually changed - so nothing\n// further needs to be typed there.\n\nconst body = context.document.body;\n\nconst dateResults = body.search(\"30/02/2025\", { matchCase: true, matchWholeWord: false });\ndateResults.load(\"text\");\nawait context.sync();\n\nif (dateResults.items.length === 0) {\n  throw new Error(\"Could not find the date '30/02/2025' to fix.\");\n}\n\n// Replace just the wrong text in place; keeps the surrounding run's\n// formatting (pt-BR language run properties) intact.\ndateResults.items[0].insertText(\"30/03/2025\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Fix the typo'd daily-standup date 30/02/2025 -> 30/03/2025 in the\n# \"ATA 30-03-2025\" minutes document (see commit \"Atualizando ata de daily\").\n# The participant name runs (\"Ana Karoline, Lays Abreu, Vitor Restini\" and\n# \"Ana Beatriz Zinatto, Luiz Felipe\") keep the exact same text in the\n# source diff - only Word's live spell-checker re-flagged the proper\n# nouns (wrapping them in proofErr markers and re-splitting the runs\n# around them) while resaving, no characters actually changed - so\n# nothing further needs to be typed there.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"30/02/2025\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"30/03/2025\"\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceAll = 2\n$found = $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n\nif (-not $found) {\n    throw \"Could not find the date '30/02/2025' to fix.\"\n}\n"}
